# GradesDatabase: add a 4th project/assignment grading column to the
# IndividualContribs and IndividualGrades sheets.
#
# Order of operations mirrors the original commit: IndividualContribs
# ("PROJECT 4") was edited first, then IndividualGrades ("ASSIGNMENT 4") -
# this is visible in the shared-string table ordering (PROJECT 4 before
# ASSIGNMENT 4) and in the style table (the bold Calibri 11 header font
# created for IndividualContribs precedes the one for IndividualGrades).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) IndividualContribs ("Teams" contribution sheet) - new column "PROJECT 4"
# ---------------------------------------------------------------------
$wsContribs = $wb.Worksheets.Item("IndividualContribs")

$wsContribs.Range("E1").Value = "PROJECT 4"
$wsContribs.Range("E1").Font.Bold = $true
$wsContribs.Range("E1").Font.Name = "Calibri"
$wsContribs.Range("E1").Font.Size = 11

# Rows 2-6 only go up to column C so far, the new grade lands in column D.
$wsContribs.Range("D2").Value = 95
$wsContribs.Range("D3").Value = 96
$wsContribs.Range("D4").Value = 65
$wsContribs.Range("D5").Value = 17
$wsContribs.Range("D6").Value = 87

# Rows 7-10 already have a value through column D, so the new grade lands
# in column E.
$wsContribs.Range("E7").Value = 65
$wsContribs.Range("E8").Value = 90
$wsContribs.Range("E9").Value = 35
$wsContribs.Range("E10").Value = 67

# ---------------------------------------------------------------------
# 2) IndividualGrades - new column "ASSIGNMENT 4"
# ---------------------------------------------------------------------
$wsGrades = $wb.Worksheets.Item("IndividualGrades")

$wsGrades.Range("E1").Value = "ASSIGNMENT 4"
$wsGrades.Range("E1").Font.Bold = $true
$wsGrades.Range("E1").Font.Name = "Calibri"
$wsGrades.Range("E1").Font.Size = 11

$wsGrades.Range("E2").Value = 95
$wsGrades.Range("E3").Value = 96
$wsGrades.Range("E4").Value = 65
$wsGrades.Range("E5").Value = 17
$wsGrades.Range("E6").Value = 87
$wsGrades.Range("E7").Value = 65
$wsGrades.Range("E8").Value = 90
$wsGrades.Range("E9").Value = 35
$wsGrades.Range("E10").Value = 67
$wsGrades.Range("E11").Value = 95
$wsGrades.Range("E12").Value = 96
$wsGrades.Range("E13").Value = 65
$wsGrades.Range("E14").Value = 17
$wsGrades.Range("E15").Value = 87
